$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.848.56"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$ws.Range("D3").Value = "2.498.19"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.47%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.00%  "

# Row 9
$ws.Range("E9").Value = "  +0.91%  "

# Row 10
$ws.Range("E10").Value = "  -1.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("D11").ClearFormats()

# Row 12
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("D13").Value = "2.938.86"
$ws.Range("E13").Value = "  -0.20%  "

# Row 14
$ws.Range("D14").Value = "58.764.52"
$ws.Range("E14").Value = "  +0.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.70"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.40%  "

# Row 16
$ws.Range("E16").Value = "  -0.15%  "

# Row 17
$ws.Range("D17").Value = "2.492.90"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("E19").Value = "  +0.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("E22").Value = "  +1.66%  "

# Row 23
$ws.Range("E23").Value = "  +1.94%  "

# Row 24
$ws.Range("E24").Value = "  +0.90%  "

# Row 25
$ws.Range("E25").Value = "  -0.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.40%  "

# Row 27
$ws.Range("E27").Value = "  -0.54%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0759"
$ws.Range("E28").Value = "  -0.91%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.44"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.84%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.52%  "

# Row 31
$ws.Range("E31").Value = "  -1.03%  "

# Row 32
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.56%  "

# Row 35
$ws.Range("E35").Value = "  -1.13%  "

# Row 36
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("E37").Value = "  -2.57%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.798"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "280.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.15%  "

# Row 41
$ws.Range("E41").Value = "  +0.34%  "

# Row 42
$ws.Range("E42").Value = "  -3.48%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.29%  "

# Row 45
$ws.Range("E45").Value = "  +0.34%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0926"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
$ws.Range("E47").Value = "  -2.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0218"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("D50").Value = "1.748.16"
$ws.Range("E50").Value = "  -1.00%  "

# Row 51
$ws.Range("E51").Value = "  -0.41%  "
